$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty column C so Longitud/Tiempo/etc. shift left (matches the
# author's column layout change: D->C, E->D, F->E, G->F, H->G, I->H, J->I, K->J)
$ws.Range("C1").EntireColumn.Delete()

# --- Row 2 (Tierra): flag this row as the best match -> "SATURNO"
$ws.Range("J2").Value = "SATURNO"

# --- Row 4 (Jupiter) new data/formulas
$ws.Range("C4").Value = 3.5
$ws.Range("D4").Formula = "= C4*0.46"
$ws.Range("E4").Formula = "=(C4*2)/D4^2"
$ws.Range("F4").Formula = "=ABS(E4-9.8)/9.8"
$ws.Range("G4").Formula = "=ABS(E4-9)/9"
$ws.Range("H4").Formula = "=ABS(E4-8.7)/8.7"
$ws.Range("I4").Formula = "=ABS(E4-11)/11"

# --- Row 5 (Saturno) new data/formulas
$ws.Range("C5").Value = 2.5
$ws.Range("D5").Formula = "= C5*0.46"
$ws.Range("E5").Formula = "=(C5*2)/D5^2"
$ws.Range("F5").Formula = "=ABS(E5-9.8)/9.8"
$ws.Range("G5").Formula = "=ABS(E5-9)/9"
$ws.Range("H5").Formula = "=ABS(E5-8.7)/8.7"
$ws.Range("I5").Formula = "=ABS(E5-11)/11"

# --- Row 6 (Urano) new data/formulas
$ws.Range("C6").Value = 3
$ws.Range("D6").Formula = "= C6*0.46"
$ws.Range("E6").Formula = "=(C6*2)/D6^2"
$ws.Range("F6").Formula = "=ABS(E6-9.8)/9.8"
$ws.Range("G6").Formula = "=ABS(E6-9)/9"
$ws.Range("H6").Formula = "=ABS(E6-8.7)/8.7"
$ws.Range("I6").Formula = "=ABS(E6-11)/11"

# --- Row 7 (Neptuno) new data/formulas
$ws.Range("C7").Value = 4
$ws.Range("D7").Formula = "=(C7*D2)"
$ws.Range("E7").Formula = "=(C7*2)/D7^2"
$ws.Range("F7").Formula = "=ABS(E7-B2)/B2"
$ws.Range("G7").Formula = "=ABS(E7-9)/9"
$ws.Range("H7").Formula = "=ABS(E7-8.7)/8.7"
$ws.Range("I7").Formula = "=ABS(E7-11)/11"

# --- Row 9: totals
$ws.Range("E9").Value = "Error total"
$ws.Range("F9").Formula = "=F2+F3+F4+F5+F6+F7"
$ws.Range("G9").Formula = "=G2+G3+G4+G5+G6+G7"
$ws.Range("H9").Formula = "=H2+H3+H4+H5+H6+H7"
$ws.Range("I9").Formula = "=I2+I3+I4+I5+I6+I7"

# Selection moved to K1 in the saved file
$ws.Range("K1").Select() | Out-Null
